# Commit: "Changed Lead and Owner from excel"
#
# On the "PerformanceTC" sheet:
#  - Lead (column C, rows 2-12) was "Nadeem" -> now "Gaurav"
#  - Owner (column D, row 2) was "Sanket" -> now "Gaurav" (matching the rest of the column)
#  - CaseToRun (column B, row 2) was "Y" -> now "N"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PerformanceTC")

# Row 2: flip CaseToRun to "N" and align Owner with the rest of the sheet.
$ws.Range("B2").Value = "N"
$ws.Range("D2").Value = "Gaurav"

# Lead changes from "Nadeem" to "Gaurav" for every data row (2-12).
$ws.Range("C2:C12").Value = "Gaurav"

# Leave the cursor where the author last clicked while making the edit.
$ws.Range("E11").Select() | Out-Null
